# CodigosCuentasEstaciones.xlsx - fix: cambios para que cojieran base 21
# negativo e iva21 negativo.
#
# The "CuentasEstaciones" table (sheet "Hoja1") maps a numeric "Codigo" to
# an accounting "Cuenta" / "Empresa" pair. This edit adds a brand-new
# lookup row (Codigo 99685, Cuenta 60002995, Empresa E29) using the
# previously-unused row 23 of the table, and re-sequences the "Cuenta"
# values of rows 17-22 by one slot so the newly freed account code can be
# reused at the bottom of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sequence the existing "Cuenta" column for rows 17-22 (each one takes
# over the value previously held by the row below it).
$ws.Range("B17").Value2 = "60002595"
$ws.Range("B18").Value2 = "60002695"
$ws.Range("B19").Value2 = "60002795"
$ws.Range("B20").Value2 = "60002895"
$ws.Range("B21").Value2 = "60002995"
$ws.Range("B22").Value2 = "60003195"

# Populate the new row 23: Codigo / Cuenta / Empresa.
$ws.Range("A23").Value2 = "99685"
$ws.Range("B23").Value2 = "60002995"
$ws.Range("C23").Value2 = "E29"

# Leave the selection on the newly completed row, matching the saved
# workbook's cursor position.
$null = $ws.Range("A23").Select()
